# Update the division drill worksheet: replace the filled-in problems with
# a freshly generated set (same table shape, new "NN÷N=" prompts).
#
# Each old prompt is unique in the document, so a straight Find/Replace
# (not "replace all", just the single first match) per pair is safe and
# unambiguous, even though a couple of the new values happen to equal some
# other cell's *old* value (e.g. "10÷7=" and "60÷2=" reappear as targets).
# Processing the pairs once, in document order, guarantees we never touch
# a cell we already rewrote.

$d = $word.ActiveDocument

$replacements = @(
    @("61÷6=", "65÷8="),
    @("81÷6=", "95÷5="),
    @("10÷7=", "80÷9="),
    @("60÷2=", "80÷6="),
    @("57÷4=", "10÷7="),
    @("40÷6=", "41÷5="),
    @("75÷8=", "71÷7="),
    @("61÷5=", "70÷9="),
    @("54÷9=", "10÷2="),
    @("72÷5=", "60÷2="),
    @("56÷7=", "94÷4="),
    @("83÷5=", "35÷9="),
    @("29÷2=", "46÷4="),
    @("86÷7=", "38÷7="),
    @("45÷2=", "65÷2="),
    @("24÷4=", "90÷3="),
    @("57÷9=", "45÷7="),
    @("33÷8=", "39÷7="),
    @("66÷6=", "14÷9="),
    @("85÷9=", "93÷9="),
    @("87÷6=", "17÷2="),
    @("70÷8=", "39÷4="),
    @("52÷5=", "36÷8="),
    @("69÷2=", "68÷3="),
    @("97÷4=", "88÷6=")
)

foreach ($pair in $replacements) {
    $oldText = $pair[0]
    $newText = $pair[1]
    $found = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        Write-Host "WARNING: could not find '$oldText' to replace with '$newText'"
    }
}

Write-Host "Done: applied $($replacements.Count) replacements."
